$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$ws.Cells.Item($row, 1).Value = "Record"
$ws.Cells.Item($row, 2).Value = "RJ Record"
$ws.Cells.Item($row, 3).Value = "Social"
$ws.Cells.Item($row, 4).Value = "2025-04-02T19:08"
$ws.Cells.Item($row, 5).Value = "Positivo"
$ws.Cells.Item($row, 6).Value = "Prefeitura de Campos anuncia liberação do saldo do cartão Goitacá. Programa de transferência de renda no valor de R`$ 200. Atualmente, cerca de 20 mil famílias são atendidas pelo programa. *nota coberta*"
